$d = $word.ActiveDocument

# The <id>...</id> tag for p106v_1 is currently split across three runs
# (the literal "<id>" tag, the bare id text, and the literal "</id>" tag)
# because the id text itself has different run formatting than the
# surrounding tag markup. Collapse them into a single run carrying the
# tag-markup formatting (Courier New / color 7f6000 / sz 18) by doing a
# plain Find & Replace across the whole "<id>p106v_1</id>" span: Word
# merges matched text that spans multiple runs into one run using the
# formatting of the first run in the match.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("<id>p106v_1</id>", $false, $false, $false, $false, $false, `
              $true, 1, $false, "<id>p106v_1</id>", 2)
